$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Turns Played") values: set B11:F11 to 0
$ws.Range("B11:F11").Value = 0

# Update the selected cell to D11
$ws.Range("D11").Select()
